# Edit applies report fixes:
#  - "Đơn 1 bác sĩ": add a new order row (HD-LUXURY 633) and push the
#    "Tổng" (total) row down, recomputing its totals.
#  - "Đơn thu nợ": add four new debt-collection rows (TN 177-180) and
#    push the "Tổng" row down, recomputing its totals.
#  - "Lương": refresh computed payroll figures to reflect the new orders.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: "Đơn 1 bác sĩ"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Đơn 1 bác sĩ")

# New row 4: HD-LUXURY 633
$ws.Cells.Item(4, 1).Value = "HD-LUXURY"
$ws.Cells.Item(4, 2).Value = 633
$ws.Cells.Item(4, 3).NumberFormat = "@"
$ws.Cells.Item(4, 3).Value = "08-07-2024"
$ws.Cells.Item(4, 4).Value = "CẦN THƠ"
$ws.Cells.Item(4, 5).Value = "Bạch Nhi"
$ws.Cells.Item(4, 6).Value = "Cá nhân"
$ws.Cells.Item(4, 7).Value = "Nâng mũi"
$ws.Cells.Item(4, 8).Value = 15000000
$ws.Cells.Item(4, 9).Value = $null
$ws.Cells.Item(4, 10).Value = $null
$ws.Cells.Item(4, 11).Value = 15000000
$ws.Cells.Item(4, 12).Value = 7000000
$ws.Cells.Item(4, 13).Value = 0.1
$ws.Cells.Item(4, 14).Value = 700000

# Row 5 now holds the recomputed "Tổng" (total) row
$ws.Cells.Item(5, 1).Value = "Tổng"
$ws.Cells.Item(5, 2).Value = 3
$ws.Cells.Item(5, 3).Value = $null
$ws.Cells.Item(5, 4).Value = $null
$ws.Cells.Item(5, 5).Value = $null
$ws.Cells.Item(5, 6).Value = $null
$ws.Cells.Item(5, 7).Value = $null
$ws.Cells.Item(5, 8).Value = 52000000
$ws.Cells.Item(5, 9).Value = $null
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 52000000
$ws.Cells.Item(5, 12).Value = 44000000
$ws.Cells.Item(5, 13).Value = 0
$ws.Cells.Item(5, 14).Value = 4400000

# ---------------------------------------------------------------------
# Sheet: "Đơn thu nợ"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Đơn thu nợ")

# New row 3: TN 177
$ws.Cells.Item(3, 1).Value = "TN"
$ws.Cells.Item(3, 2).Value = 177
$ws.Cells.Item(3, 3).Value = 8000000
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "08-05-2024"
$ws.Cells.Item(3, 5).Value = "SÓC TRĂNG"
$ws.Cells.Item(3, 6).Value = "HD-LUXURY-611"
$ws.Cells.Item(3, 7).Value = "Nâng mũi"
$ws.Cells.Item(3, 8).Value = "ngọc hân"
$ws.Cells.Item(3, 9).Value = "CTV"
$ws.Cells.Item(3, 10).Value = "Thạch Hoàng Nhân"
$ws.Cells.Item(3, 11).Value = 35000000
$ws.Cells.Item(3, 12).Value = "Lê Đình Hậu"
$ws.Cells.Item(3, 13).Value = 8000000
$ws.Cells.Item(3, 14).Value = 43000000
$ws.Cells.Item(3, 15).Value = 43000000
$ws.Cells.Item(3, 16).Value = "Phạm Thanh Hoàng"
$ws.Cells.Item(3, 17).Value = $null
$ws.Cells.Item(3, 18).Value = 0
$ws.Cells.Item(3, 19).Value = 0
$ws.Cells.Item(3, 20).Value = 0
$ws.Cells.Item(3, 21).Value = 0
$ws.Cells.Item(3, 22).Value = 0.1
$ws.Cells.Item(3, 23).Value = 800000
$ws.Cells.Item(3, 24).Value = 0
$ws.Cells.Item(3, 25).Value = 0

# New row 4: TN 178
$ws.Cells.Item(4, 1).Value = "TN"
$ws.Cells.Item(4, 2).Value = 178
$ws.Cells.Item(4, 3).Value = 2000000
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "08-06-2024"
$ws.Cells.Item(4, 5).Value = "SÓC TRĂNG"
$ws.Cells.Item(4, 6).Value = "HD-LUXURY-356"
$ws.Cells.Item(4, 7).Value = "Nâng mũi"
$ws.Cells.Item(4, 8).Value = "thuý vân"
$ws.Cells.Item(4, 9).Value = "CTV"
$ws.Cells.Item(4, 10).Value = "Thạch Hoàng Nhân"
$ws.Cells.Item(4, 11).Value = 35000000
$ws.Cells.Item(4, 12).Value = $null
$ws.Cells.Item(4, 13).Value = $null
$ws.Cells.Item(4, 14).Value = 35000000
$ws.Cells.Item(4, 15).Value = 26000000
$ws.Cells.Item(4, 16).Value = "Phạm Thanh Hoàng"
$ws.Cells.Item(4, 17).Value = $null
$ws.Cells.Item(4, 18).Value = 0
$ws.Cells.Item(4, 19).Value = 0
$ws.Cells.Item(4, 20).Value = 0
$ws.Cells.Item(4, 21).Value = 0
$ws.Cells.Item(4, 22).Value = 0.1
$ws.Cells.Item(4, 23).Value = 200000
$ws.Cells.Item(4, 24).Value = 0
$ws.Cells.Item(4, 25).Value = 0

# New row 5: TN 179
$ws.Cells.Item(5, 1).Value = "TN"
$ws.Cells.Item(5, 2).Value = 179
$ws.Cells.Item(5, 3).Value = $null
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "08-09-2024"
$ws.Cells.Item(5, 5).Value = "SÓC TRĂNG"
$ws.Cells.Item(5, 6).Value = "HD-LUXURY-500"
$ws.Cells.Item(5, 7).Value = "Nâng mũi"
$ws.Cells.Item(5, 8).Value = "pola"
$ws.Cells.Item(5, 9).Value = "CTV"
$ws.Cells.Item(5, 10).Value = "Thạch Hoàng Nhân"
$ws.Cells.Item(5, 11).Value = 43000000
$ws.Cells.Item(5, 12).Value = $null
$ws.Cells.Item(5, 13).Value = $null
$ws.Cells.Item(5, 14).Value = 43000000
$ws.Cells.Item(5, 15).Value = 24000000
$ws.Cells.Item(5, 16).Value = "Phạm Thanh Hoàng"
$ws.Cells.Item(5, 17).Value = $null
$ws.Cells.Item(5, 18).Value = 0
$ws.Cells.Item(5, 19).Value = 0
$ws.Cells.Item(5, 20).Value = 0
$ws.Cells.Item(5, 21).Value = 0
$ws.Cells.Item(5, 22).Value = 0.1
$ws.Cells.Item(5, 23).Value = $null
$ws.Cells.Item(5, 24).Value = 0
$ws.Cells.Item(5, 25).Value = 0

# New row 6: TN 180
$ws.Cells.Item(6, 1).Value = "TN"
$ws.Cells.Item(6, 2).Value = 180
$ws.Cells.Item(6, 3).Value = 2000000
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "08-09-2024"
$ws.Cells.Item(6, 5).Value = "SÓC TRĂNG"
$ws.Cells.Item(6, 6).Value = "HD-LUXURY-500"
$ws.Cells.Item(6, 7).Value = "Nâng mũi"
$ws.Cells.Item(6, 8).Value = "pola"
$ws.Cells.Item(6, 9).Value = "CTV"
$ws.Cells.Item(6, 10).Value = "Thạch Hoàng Nhân"
$ws.Cells.Item(6, 11).Value = 43000000
$ws.Cells.Item(6, 12).Value = $null
$ws.Cells.Item(6, 13).Value = $null
$ws.Cells.Item(6, 14).Value = 43000000
$ws.Cells.Item(6, 15).Value = 24000000
$ws.Cells.Item(6, 16).Value = "Phạm Thanh Hoàng"
$ws.Cells.Item(6, 17).Value = $null
$ws.Cells.Item(6, 18).Value = 0
$ws.Cells.Item(6, 19).Value = 0
$ws.Cells.Item(6, 20).Value = 0
$ws.Cells.Item(6, 21).Value = 0
$ws.Cells.Item(6, 22).Value = 0.1
$ws.Cells.Item(6, 23).Value = 200000
$ws.Cells.Item(6, 24).Value = 0
$ws.Cells.Item(6, 25).Value = 0

# Row 7 now holds the recomputed "Tổng" (total) row
$ws.Cells.Item(7, 1).Value = "Tổng"
$ws.Cells.Item(7, 2).Value = 5
$ws.Cells.Item(7, 3).Value = 22000000
$ws.Cells.Item(7, 4).Value = $null
$ws.Cells.Item(7, 5).Value = $null
$ws.Cells.Item(7, 6).Value = $null
$ws.Cells.Item(7, 7).Value = $null
$ws.Cells.Item(7, 8).Value = $null
$ws.Cells.Item(7, 9).Value = $null
$ws.Cells.Item(7, 10).Value = $null
$ws.Cells.Item(7, 11).Value = 176000000
$ws.Cells.Item(7, 12).Value = $null
$ws.Cells.Item(7, 13).Value = 8000000
$ws.Cells.Item(7, 14).Value = 184000000
$ws.Cells.Item(7, 15).Value = 137000000
$ws.Cells.Item(7, 16).Value = $null
$ws.Cells.Item(7, 17).Value = $null
$ws.Cells.Item(7, 18).Value = 0
$ws.Cells.Item(7, 19).Value = 0
$ws.Cells.Item(7, 20).Value = 0
$ws.Cells.Item(7, 21).Value = 0
$ws.Cells.Item(7, 22).Value = 0
$ws.Cells.Item(7, 23).Value = 2200000
$ws.Cells.Item(7, 24).Value = 0
$ws.Cells.Item(7, 25).Value = 0

# ---------------------------------------------------------------------
# Sheet: "Lương" (refresh computed payroll totals)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Lương")

$ws.Cells.Item(2, 2).Value = 4
$ws.Cells.Item(3, 2).Value = 140000
$ws.Cells.Item(4, 2).Value = 428571.4285714285
$ws.Cells.Item(7, 2).Value = 4400000
$ws.Cells.Item(12, 2).Value = -2000000
$ws.Cells.Item(15, 2).Value = 428571.4285714285
$ws.Cells.Item(26, 2).Value = 428571.4285714285
$ws.Cells.Item(33, 2).Value = 1200000
$ws.Cells.Item(35, 2).Value = 7968571.428571429
$ws.Cells.Item(36, 2).Value = 428571.4285714285
$ws.Cells.Item(37, 2).Value = 1628571.428571429
$ws.Cells.Item(38, 2).Value = 10025714.28571429
